# Auto-generated edit script: apply "Add data for 2024-01-25" crime-data refresh
# Updates 155 cells across 42 worksheets (2023 corrections + new 2024 YTD counts).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7700
$ws.Range("K2").Value = 413
$ws.Range("J3").Value = 8071
$ws.Range("K3").Value = 383
$ws.Range("K4").Value = 75
$ws.Range("J6").Value = 11051
$ws.Range("K6").Value = 523
$ws.Range("J7").Value = 29219
$ws.Range("K7").Value = 1415

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 14
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 39
$ws.Range("K8").Value = 92
$ws.Range("K9").Value = 5
$ws.Range("K14").Value = 8
$ws.Range("J15").Value = 360
$ws.Range("K15").Value = 7
$ws.Range("K19").Value = 30
$ws.Range("K20").Value = 38
$ws.Range("K23").Value = 13
$ws.Range("K24").Value = 5
$ws.Range("K29").Value = 74
$ws.Range("K31").Value = 20
$ws.Range("K33").Value = 67
$ws.Range("K37").Value = 38
$ws.Range("K41").Value = 15
$ws.Range("K42").Value = 50
$ws.Range("J43").Value = 246
$ws.Range("K43").Value = 16
$ws.Range("K44").Value = 11
$ws.Range("K48").Value = 15
$ws.Range("K50").Value = 6
$ws.Range("K54").Value = 23
$ws.Range("K56").Value = 2
$ws.Range("J63").Value = 85
$ws.Range("K63").Value = 10
$ws.Range("K65").Value = 40
$ws.Range("K67").Value = 53
$ws.Range("K69").Value = 3
$ws.Range("K72").Value = 5
$ws.Range("K75").Value = 5
$ws.Range("J79").Value = 800
$ws.Range("K79").Value = 31
$ws.Range("K80").Value = 3
$ws.Range("K82").Value = 2
$ws.Range("K85").Value = 68
$ws.Range("K88").Value = 24
$ws.Range("K89").Value = 22
$ws.Range("K90").Value = 11
$ws.Range("K91").Value = 15
$ws.Range("K92").Value = 6
$ws.Range("K95").Value = 28
$ws.Range("K96").Value = 20
$ws.Range("K97").Value = 12
$ws.Range("J101").Value = 29219
$ws.Range("K101").Value = 1415

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 8

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 3

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 25
$ws.Range("K3").Value = 29
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 13
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 211
$ws.Range("J3").Value = 193
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 234
$ws.Range("K3").Value = 10
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J4").Value = 15
$ws.Range("K6").Value = 6
$ws.Range("J7").Value = 360
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 2
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 147
$ws.Range("K6").Value = 7
$ws.Range("J7").Value = 246
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 1
$ws.Range("K7").Value = 5

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 2

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 2

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 3
